# chore: update Sheets via scheduled runner
# Applies updated market/profit figures to several leve rows across the
# job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

function Set-Cells {
    param(
        $ws,
        [int]$row,
        [hashtable]$values
    )
    foreach ($col in $values.Keys) {
        $ws.Cells.Item($row, $col).Value = $values[$col]
    }
}

# Column map: H=8 I=9 J=10 K=11 L=12 M=13 N=14

# ---------------- ALC ----------------
$ws = $wb.Worksheets.Item("ALC")

Set-Cells $ws 15 @{ 8=1727.4375; 9=1727.4375; 11=5182.3125; 13=-5013.3125 }
Set-Cells $ws 74 @{ 8=3920763.5; 9=4247077; 10=5000; 11=4247077; 12=5000; 13=-4246141; 14=-6872 }
Set-Cells $ws 76 @{ 8=157145550; 9=366667840; 10=3848.5; 11=366667840; 12=3848.5; 13=-366667525; 14=-4478.5 }
Set-Cells $ws 77 @{ 8=3920763.5; 9=4247077; 10=5000; 11=21235385; 12=25000; 13=-21230705; 14=-34360 }
Set-Cells $ws 79 @{ 8=157145550; 9=366667840; 10=3848.5; 11=366667840; 12=3848.5; 13=-366666748; 14=-6032.5 }
Set-Cells $ws 129 @{ 8=733.8182; 10=1540; 12=4620; 14=-14620 }
Set-Cells $ws 137 @{ 8=18520108; 9=1637.875; 10=166667870; 11=4913.625; 12=500003610; 13=-2363.625; 14=-500008710 }
Set-Cells $ws 138 @{ 8=4575.89; 9=1847.5; 10=4689.5728; 11=5542.5; 12=14068.7184; 13=-402.5; 14=-24348.7184 }

# ---------------- ARM ----------------
$ws = $wb.Worksheets.Item("ARM")

Set-Cells $ws 32 @{ 8=18179.188; 9=18466.111; 11=18466.111; 13=-18179.111 }
Set-Cells $ws 61 @{ 8=1331.3334; 9=1183.875; 11=1183.875; 13=-971.875 }
Set-Cells $ws 74 @{ 8=771.1539; 9=758.44446; 10=799.75; 11=758.44446; 12=799.75; 13=115.55554; 14=-2547.75 }
Set-Cells $ws 77 @{ 8=771.1539; 9=758.44446; 10=799.75; 11=3792.2223; 12=3998.75; 13=575.7776999999996; 14=-12734.75 }
Set-Cells $ws 136 @{ 8=1331.3334; 9=1183.875; 11=3551.625; 13=-1001.625 }

# ---------------- BSM ----------------
$ws = $wb.Worksheets.Item("BSM")

Set-Cells $ws 22 @{ 8=100; 9=100; 11=100; 13=73 }
Set-Cells $ws 105 @{ 8=2804.05; 9=2746.8; 10=2975.8; 11=2746.8; 12=2975.8; 13=-999.8000000000002; 14=-6469.8 }
Set-Cells $ws 107 @{ 8=1272.8462; 9=1119.7; 11=1119.7; 13=800.3 }

# ---------------- CRP ----------------
$ws = $wb.Worksheets.Item("CRP")

Set-Cells $ws 58 @{ 8=3811.9143; 9=978.25; 10=9994.454; 11=978.25; 12=9994.454; 13=-775.25; 14=-10400.454 }
Set-Cells $ws 62 @{ 8=7358.4346; 9=8091.923; 10=6404.9; 11=8091.923; 12=6404.9; 13=-7467.923; 14=-7652.9 }
Set-Cells $ws 65 @{ 8=7358.4346; 9=8091.923; 10=6404.9; 11=40459.615; 12=32024.5; 13=-37339.615; 14=-38264.5 }
Set-Cells $ws 70 @{ 8=13800; 10=13800; 12=13800; 14=-14430 }
Set-Cells $ws 73 @{ 8=13800; 10=13800; 12=13800; 14=-15984 }
Set-Cells $ws 86 @{ 8=2775.5; 9=3058.55; 11=3058.55; 13=-1935.55 }
Set-Cells $ws 89 @{ 8=2775.5; 9=3058.55; 11=15292.75; 13=-9676.75 }
Set-Cells $ws 99 @{ 8=2165.3845; 9=2112.5; 11=2112.5; 13=-614.5 }
Set-Cells $ws 126 @{ 8=2165.3845; 9=2112.5; 11=6337.5; 13=-3867.5 }
Set-Cells $ws 134 @{ 8=1546.0344; 9=1345.1364; 11=4035.4092; 13=-1500.4092 }
Set-Cells $ws 136 @{ 8=3811.9143; 9=978.25; 10=9994.454; 11=2934.75; 12=29983.362; 13=-384.75; 14=-35083.362 }

# ---------------- CUL ----------------
$ws = $wb.Worksheets.Item("CUL")

Set-Cells $ws 131 @{ 8=13193.07; 10=1357.0853; 12=4071.2559; 14=-14151.2559 }

# ---------------- GSM ----------------
$ws = $wb.Worksheets.Item("GSM")

Set-Cells $ws 70 @{ 8=4966.6665; 9=0; 11=0 }
$ws.Cells.Item(70, 13).ClearContents()
Set-Cells $ws 73 @{ 8=4966.6665; 9=0; 11=0 }
$ws.Cells.Item(73, 13).ClearContents()
Set-Cells $ws 113 @{ 8=1014.8; 9=923.25; 10=1177.5555; 11=923.25; 12=1177.5555; 13=1246.75; 14=-5517.5555 }

# ---------------- LTW ----------------
$ws = $wb.Worksheets.Item("LTW")

Set-Cells $ws 68 @{ 8=2100.3333; 9=2040.4; 10=2400; 11=2040.4; 12=2400; 13=-1291.4; 14=-3898 }
Set-Cells $ws 71 @{ 8=2100.3333; 9=2040.4; 10=2400; 11=10202; 12=12000; 13=-6458; 14=-19488 }
Set-Cells $ws 132 @{ 8=4417.5884; 9=4859.6; 11=14578.8; 13=-12048.8 }

# ---------------- WVR ----------------
$ws = $wb.Worksheets.Item("WVR")

Set-Cells $ws 75 @{ 8=20000; 10=20000; 12=20000; 14=-21872 }
Set-Cells $ws 78 @{ 8=20000; 10=20000; 12=60000; 14=-69360 }
Set-Cells $ws 113 @{ 8=461.92307; 9=488; 10=375; 11=1464; 12=1125; 13=706; 14=-5465 }

Write-Host "Applied Bahamut_Profits leve-profit updates."
